$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New 2016 observations, one per district, in the same order as the existing
# 2009 rows on the "Data" sheet (rows 2..48).
# Columns: district_id, district_name, 2016 value
# ---------------------------------------------------------------------------
$districts2016 = @(
    @('d18974', 'Baringo', 3.2),
    @('d18975', 'Bomet', 0.7),
    @('d18987', 'Bungoma', 2.1),
    @('d18988', 'Busia', 3.3),
    @('d18976', 'Elgeyo-Marakwet', 2),
    @('d18955', 'Embu', 4),
    @('d18965', 'Garissa', 0.9),
    @('d18968', 'Homa Bay', 2.9),
    @('d18956', 'Isiolo', 1.7),
    @('d18991', 'Kajiado', 0.7),
    @('d18989', 'Kakamega', 2.1),
    @('d18977', 'Kericho', 2.1),
    @('d18943', 'Kiambu', 1.1),
    @('d18949', 'Kilifi', 2.7),
    @('d18944', 'Kirinyaga', 5.8),
    @('d18969', 'Kisii', 2.4),
    @('d18970', 'Kisumu', 2.7),
    @('d18957', 'Kitui', 4.6),
    @('d18950', 'Kwale', 4.1),
    @('d18978', 'Laikipia', 1.6),
    @('d18951', 'Lamu', 2.9),
    @('d18958', 'Machakos', 3.9),
    @('d18959', 'Makueni', 3.9),
    @('d18966', 'Mandera', 3.7),
    @('d18960', 'Marsabit', 3.6),
    @('d18961', 'Meru', 8.5),
    @('d18971', 'Migori', 5.1),
    @('d18952', 'Mombasa', 0.5),
    @('d18946', 'Murang''a', 3.3),
    @('d18964', 'Nairobi', 1.2),
    @('d18979', 'Nakuru', 1.9),
    @('d18980', 'Nandi', 1.1),
    @('d18981', 'Narok', 1.9),
    @('d18962', 'Nithi', 4.4),
    @('d18972', 'Nyamira', 5.2),
    @('d18947', 'Nyandarua', 1.7),
    @('d18948', 'Nyeri', 2.1),
    @('d18982', 'Samburu', 1.9),
    @('d18973', 'Siaya', 5.3),
    @('d18953', 'Taita Taveta', 3),
    @('d18954', 'Tana River', 2.8),
    @('d18983', 'Trans-Nzoia', 2.1),
    @('d18984', 'Turkana', 4),
    @('d18985', 'Uasin Gishu', 3.6),
    @('d18990', 'Vihiga', 6.3),
    @('d18967', 'Wajir', 1.5),
    @('d18986', 'West Pokot', 2.2)

)

$dataSheet = $wb.Worksheets.Item("Data")

# Insert a blank row directly below each existing district row, working from
# the bottom of the table upwards so the row numbers we still need to touch
# never shift underneath us.
$firstRow = 2
$lastRow = 48
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $dataSheet.Rows.Item($r + 1).Insert()
}

# Fill in the freshly inserted rows (3, 5, 7, ... 95) with the 2016 data.
for ($i = 0; $i -lt $districts2016.Count; $i++) {
    $targetRow = $firstRow + 1 + (2 * $i)
    $entry = $districts2016[$i]
    $dataSheet.Cells.Item($targetRow, 1).Value = $entry[0]
    $dataSheet.Cells.Item($targetRow, 2).Value = $entry[1]
    $dataSheet.Cells.Item($targetRow, 3).Value = 2016
    $dataSheet.Cells.Item($targetRow, 4).Value = $entry[2]
}

# ---------------------------------------------------------------------------
# "Data-wide-value" sheet: add a "2016" column (C) next to the existing
# "2009" value column (B), one value per district in the same row order.
# The header "2016" must stay a text label, not become the number 2016, so
# the cell is briefly formatted as Text while the label is entered and the
# formatting is cleared again afterwards (leaving the cell on the default
# style but with a genuine string value).
# ---------------------------------------------------------------------------
$wideSheet = $wb.Worksheets.Item("Data-wide-value")
$headerCell = $wideSheet.Cells.Item(1, 3)
$headerCell.NumberFormat = "@"
$headerCell.Value = "2016"
$headerCell.ClearFormats()

for ($i = 0; $i -lt $districts2016.Count; $i++) {
    $entry = $districts2016[$i]
    $wideSheet.Cells.Item(2 + $i, 3).Value = $entry[2]
}

# ---------------------------------------------------------------------------
# "Notes" sheet: the source note changes from "Source: 47" to "Source: 37".
# ---------------------------------------------------------------------------
$notesSheet = $wb.Worksheets.Item("Notes")
$notesSheet.Cells.Item(4, 1).Value = "Source: 37"
